{"js": "// Applies:\n//  1) Adds a `w:ind w:left=\"708\"` indent to the five \"UC\" progress\n//     paragraphs (\"bueno ya he desglozado el UC1...\", \"UC2, de igual\n//     manera...\", \"El UC3 tuve un poquito...\", \"UC 4\", \"UC 5\").\n//  2) Splits the trailing run of the very last paragraph so \"este\" is\n//     wrapped in proofErr spell-check markers, and appends the new\n//     30/11/2024 journal entry (several new paragraphs, one of them\n//     empty, one indented, one with a literal tab) after it.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// ---- 1) indent the five \"UC\" paragraphs -------------------------------\nconst TWIPS_PER_POINT = 20;\nconst indentTargets = [\n  \"bueno ya he \",\n  \"UC2, de igual manera, al hacer un desglose \",\n  \"El UC3 tuve un poquito \",\n  \"UC 4\",\n  \"UC 5\",\n];\n\nfor (const item of paragraphs.items) {\n  const text = item.text;\n  for (const target of indentTargets) {\n    if (text.indexOf(target) === 0) {\n      item.leftIndent = 708 / TWIPS_PER_POINT;\n      break;\n    }\n  }\n}\nawait context.sync();\n\n// ---- 2) rewrite the closing paragraph + append the new entry ----------\nconst lastParagraph = paragraphs.items[paragraphs.items.length - 1];\n\nconst ooxmlOpen =\n  '<?xml version=\"1.0\"?>' +\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" ' +\n  'pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  '<pkg:xmlData>' +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  \"<w:body>\";\nconst ooxmlClose = \"</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>\";\n\n// This first <w:p> replaces the current (last) paragraph's content --\n// the split introduces the proofErr spell-check wrapper around \"este\".\nconst closingParagraph =\n  \"<w:p>\" +\n  '<w:r><w:t xml:space=\"preserve\">Me he topado con un par de inconsistencias al enumerar los casos de uso y con un punto faltante de informaci\u00f3n, ma\u00f1ana continuare con los casos de uso, hubo un peque\u00f1o descuido que mencionare ma\u00f1ana cuando </w:t></w:r>' +\n  '<w:proofErr w:type=\"spellStart\"/>' +\n  \"<w:r><w:t>este</w:t></w:r>\" +\n  '<w:proofErr w:type=\"spellEnd\"/>' +\n  '<w:r><w:t xml:space=\"preserve\"> resuelto</w:t></w:r>' +\n  \"</w:p>\";\n\n// Remaining brand-new paragraphs for the 30/11/2024 entry.\nconst newParagraphs =\n  \"<w:p/>\" +\n  \"<w:p><w:r><w:lastRenderedPageBreak/><w:t>30/11/2024</w:t></w:r></w:p>\" +\n  \"<w:p>\" +\n  '<w:r><w:t xml:space=\"preserve\">Ok, algo paso ayer, me percate que al momento de enumerar los casos de uso, me di cuenta de que deje la </w:t></w:r>' +\n  '<w:proofErr w:type=\"spellStart\"/>' +\n  \"<w:r><w:t>evaluacion</w:t></w:r>\" +\n  '<w:proofErr w:type=\"spellEnd\"/>' +\n  '<w:r><w:t xml:space=\"preserve\"> de docente bien de lado, como si nunca hubiera existido, afortunadamente, la modificaciones a </w:t></w:r>' +\n  '<w:proofErr w:type=\"spellStart\"/>' +\n  \"<w:r><w:t>relizar</w:t></w:r>\" +\n  '<w:proofErr w:type=\"spellEnd\"/>' +\n  '<w:r><w:t xml:space=\"preserve\"> del documento no son tan extremas, de hecho, son solo algunos enumeraciones, de esto me percate ayer y fue menos dif\u00edcil cambiarlo de lo que esperaba</w:t></w:r>' +\n  \"</w:p>\" +\n  \"<w:p>\" +\n  '<w:pPr><w:ind w:left=\"708\"/></w:pPr>' +\n  \"<w:r><w:t>Que espero lograr hoy, seguir avanzando con los casos de uso del alumno y si tengo suerte, acabarlos hoy</w:t></w:r>\" +\n  \"</w:p>\" +\n  \"<w:p>\" +\n  '<w:r><w:t xml:space=\"preserve\"> </w:t></w:r>' +\n  \"<w:r><w:tab/><w:t>Corregir el tema de la evaluaci\u00f3n de docente</w:t></w:r>\" +\n  \"</w:p>\" +\n  \"<w:p>\" +\n  '<w:r><w:t xml:space=\"preserve\">De acuerdo, ya he corregido parte del tema de </w:t></w:r>' +\n  '<w:proofErr w:type=\"spellStart\"/>' +\n  \"<w:r><w:t>evaluacion</w:t></w:r>\" +\n  '<w:proofErr w:type=\"spellEnd\"/>' +\n  '<w:r><w:t xml:space=\"preserve\"> de docente junto con eso, he hecho con \u00e9xito el caso de uso relacionado al mismo del alumno, ma\u00f1ana espero continuar con los casos de uso del parte del profesor, ma\u00f1ana ser\u00e1 un </w:t></w:r>' +\n  '<w:proofErr w:type=\"spellStart\"/>' +\n  \"<w:r><w:t>dia</w:t></w:r>\" +\n  '<w:proofErr w:type=\"spellEnd\"/>' +\n  '<w:r><w:t xml:space=\"preserve\"> largo pero conque logre hacer 1 me doy por bien servido</w:t></w:r>' +\n  \"</w:p>\";\n\nlastParagraph.insertOoxml(\n  ooxmlOpen + closingParagraph + newParagraphs + ooxmlClose,\n  Word.InsertLocation.replace\n);\nawait context.sync();\n", "ps1": "# Applies:\n#  1) Adds a `w:ind w:left=\"708\"` indent to the five \"UC\" progress\n#     paragraphs (\"bueno ya he desglozado el UC1...\", \"UC2, de igual\n#     manera...\", \"El UC3 tuve un poquito...\", \"UC 4\", \"UC 5\").\n#  2) Splits the trailing run of the very last paragraph so \"este\" is\n#     wrapped in proofErr spell-check markers, and appends the new\n#     30/11/2024 journal entry (several new paragraphs, one of them\n#     empty, one indented, one with a literal tab) after it.\n\n$d = $word.ActiveDocument\n\n# ---- 1) indent the five \"UC\" paragraphs --------------------------------\nforeach ($para in $d.Paragraphs) {\n    $t = $para.Range.Text\n    if ($t.StartsWith(\"bueno ya he \") -or `\n        $t.StartsWith(\"UC2, de igual manera, al hacer un desglose \") -or `\n        $t.StartsWith(\"El UC3 tuve un poquito \") -or `\n        $t -eq \"UC 4`r\" -or `\n        $t -eq \"UC 5`r\") {\n        $para.LeftIndent = 35.4\n    }\n}\n\n# ---- 2) split \"este\" out of the closing paragraph with proofErr markup -\n$lastParagraph = $d.Paragraphs.Item($d.Paragraphs.Count)\n$fullText = $lastParagraph.Range.Text\n$target = \"este resuelto\"\n$offset = $fullText.IndexOf($target)\n$subStart = $lastParagraph.Range.Start + $offset\n$subEnd = $subStart + $target.Length\n$targetRange = $d.Range($subStart, $subEnd)\n\n$splitXml = '<?xml version=\"1.0\"?><pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\"><pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body><w:p><w:proofErr w:type=\"spellStart\"/><w:r><w:t>este</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:t xml:space=\"preserve\"> resuelto</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'\n$targetRange.InsertXML($splitXml)\n\n# ---- 3) append the new 30/11/2024 entry after the (unchanged) paragraph mark\n$lastParagraph = $d.Paragraphs.Item($d.Paragraphs.Count)\n$insertionPoint = $d.Range($lastParagraph.Range.End - 1, $lastParagraph.Range.End - 1)\n\n$newParagraphsXml = '<w:p/>' + `\n    '<w:p><w:r><w:lastRenderedPageBreak/><w:t>30/11/2024</w:t></w:r></w:p>' + `\n    '<w:p><w:r><w:t xml:space=\"preserve\">Ok, algo paso ayer, me percate que al momento de enumerar los casos de uso, me di cuenta de que deje la </w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:r><w:t>evaluacion</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:t xml:space=\"preserve\"> de docente bien de lado, como si nunca hubiera existido, afortunadamente, la modificaciones a </w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:r><w:t>relizar</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:t xml:space=\"preserve\"> del documento no son tan extremas, de hecho, son solo algunos enumeraciones, de esto me percate ayer y fue menos dif\u00edcil cambiarlo de lo que esperaba</w:t></w:r></w:p>' + `\n    '<w:p><w:pPr><w:ind w:left=\"708\"/></w:pPr><w:r><w:t>Que espero lograr hoy, seguir avanzando con los casos de uso del alumno y si tengo suerte, acabarlos hoy</w:t></w:r></w:p>' + `\n    '<w:p><w:r><w:t xml:space=\"preserve\"> </w:t></w:r><w:r><w:tab/><w:t>Corregir el tema de la evaluaci\u00f3n de docente</w:t></w:r></w:p>' + `\n    '<w:p><w:r><w:t xml:space=\"preserve\">De acuerdo, ya he corregido parte del tema de </w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:r><w:t>evaluacion</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:t xml:space=\"preserve\"> de docente junto con eso, he hecho con \u00e9xito el caso de uso relacionado al mismo del alumno, ma\u00f1ana espero continuar con los casos de uso del parte del profesor, ma\u00f1ana ser\u00e1 un </w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:r><w:t>dia</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:t xml:space=\"preserve\"> largo pero conque logre hacer 1 me doy por bien servido</w:t></w:r></w:p>'\n\n$newEntryXml = '<?xml version=\"1.0\"?><pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\"><pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body>' + $newParagraphsXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'\n$insertionPoint.InsertXML($newEntryXml)\n"}
